# Daily attendance update - 2025-08-14
# Fill in the "WFH" column (S, 2025-08-14) and the new "India Holiday"
# column (T, 2025-08-15) for every employee row (3-18) on the
# WCS_Team_August_2025 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WCS_Team_August_2025")

# Grab a cell that already carries the "Neutral" (yellow) cell style used
# for WFH/SL/Leave entries so the S-column cells pick up the same
# formatting (fill + font + border) that the rest of the sheet uses.
$ws.Range("R7").Copy()

for ($row = 3; $row -le 18; $row++) {
    $sCell = $ws.Cells.Item($row, 19)   # column S = 2025-08-14
    $tCell = $ws.Cells.Item($row, 20)   # column T = 2025-08-15 (India Holiday)

    $sCell.PasteSpecial(-4122)          # xlPasteFormats

    if ($row -eq 16) {
        $sCell.Value = "SL"
    } else {
        $sCell.Value = "WFH"
    }

    $tCell.Value = "India Holiday"
}

$excel.CutCopyMode = 0

$ws.Activate()
$ws.Range("S8").Select()
